$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 363.42856
$ws.Range("I33").Value = 316
$ws.Range("J33").Value = 648
$ws.Range("K33").Value = 316
$ws.Range("L33").Value = 648
$ws.Range("M33").Value = -87
$ws.Range("N33").Value = -1106
$ws.Range("H43").Value = 4559.8
$ws.Range("I43").Value = 1266.6666
$ws.Range("K43").Value = 1266.6666
$ws.Range("M43").Value = -1197.6666
$ws.Range("H51").Value = 6563.3335
$ws.Range("I51").Value = 6596
$ws.Range("J51").Value = 6547
$ws.Range("K51").Value = 6596
$ws.Range("L51").Value = 6547
$ws.Range("M51").Value = -6112
$ws.Range("N51").Value = -7515
$ws.Range("H52").Value = 2393.5
$ws.Range("J52").Value = 3500
$ws.Range("L52").Value = 10500
$ws.Range("N52").Value = -10820
$ws.Range("H74").Value = 10631.105
$ws.Range("I74").Value = 9317.546
$ws.Range("J74").Value = 12437.25
$ws.Range("K74").Value = 9317.546
$ws.Range("L74").Value = 12437.25
$ws.Range("M74").Value = -8381.546
$ws.Range("N74").Value = -14309.25
$ws.Range("H76").Value = 3449.5
$ws.Range("I76").Value = 3999
$ws.Range("K76").Value = 3999
$ws.Range("M76").Value = -3684
$ws.Range("H77").Value = 10631.105
$ws.Range("I77").Value = 9317.546
$ws.Range("J77").Value = 12437.25
$ws.Range("K77").Value = 46587.73
$ws.Range("L77").Value = 62186.25
$ws.Range("M77").Value = -41907.73
$ws.Range("N77").Value = -71546.25
$ws.Range("H79").Value = 3449.5
$ws.Range("I79").Value = 3999
$ws.Range("K79").Value = 3999
$ws.Range("M79").Value = -2907
$ws.Range("H113").Value = 40016092
$ws.Range("I113").Value = 100004104
$ws.Range("J113").Value = 24087.867
$ws.Range("K113").Value = 100004104
$ws.Range("L113").Value = 24087.867
$ws.Range("M113").Value = -100000850
$ws.Range("N113").Value = -30595.867
$ws.Range("H116").Value = 4428.5713
$ws.Range("I116").Value = 3200
$ws.Range("J116").Value = 5350
$ws.Range("K116").Value = 3200
$ws.Range("L116").Value = 5350
$ws.Range("M116").Value = 242
$ws.Range("N116").Value = -12234
$ws.Range("H135").Value = 1325.4615
$ws.Range("I135").Value = 976.5
$ws.Range("K135").Value = 8788.5
$ws.Range("M135").Value = -6253.5
$ws.Range("H141").Value = 4255.6665
$ws.Range("I141").Value = 3824.3635
$ws.Range("K141").Value = 11473.0905
$ws.Range("M141").Value = -6293.0905

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1892.4138
$ws.Range("I2").Value = 1625.2963
$ws.Range("K2").Value = 1625.2963
$ws.Range("M2").Value = -1512.2963
$ws.Range("H116").Value = 1892.4138
$ws.Range("I116").Value = 1625.2963
$ws.Range("K116").Value = 1625.2963
$ws.Range("M116").Value = 668.7037
$ws.Range("H122").Value = 1643.9166
$ws.Range("I122").Value = 1612
$ws.Range("K122").Value = 4836
$ws.Range("M122").Value = -2386
$ws.Range("H132").Value = 1658.6666
$ws.Range("I132").Value = 1516.6
$ws.Range("K132").Value = 4549.799999999999
$ws.Range("M132").Value = -2019.799999999999

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1892.4138
$ws.Range("I3").Value = 1625.2963
$ws.Range("K3").Value = 1625.2963
$ws.Range("M3").Value = -1511.2963

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 401.26666
$ws.Range("I22").Value = 365.63635
$ws.Range("J22").Value = 499.25
$ws.Range("K22").Value = 365.63635
$ws.Range("L22").Value = 499.25
$ws.Range("M22").Value = -15.63634999999999
$ws.Range("N22").Value = -1199.25
$ws.Range("H59").Value = 42000
$ws.Range("I59").Value = 40000
$ws.Range("K59").Value = 40000
$ws.Range("M59").Value = -38855
$ws.Range("H62").Value = 4499
$ws.Range("J62").Value = 1498
$ws.Range("L62").Value = 1498
$ws.Range("N62").Value = -2746
$ws.Range("H65").Value = 4499
$ws.Range("J65").Value = 1498
$ws.Range("L65").Value = 7490
$ws.Range("N65").Value = -13730

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 1000367.5
$ws.Range("I7").Value = 1200361.6
$ws.Range("J7").Value = 397
$ws.Range("K7").Value = 3601084.8
$ws.Range("L7").Value = 1191
$ws.Range("M7").Value = -3600972.8
$ws.Range("N7").Value = -1415
$ws.Range("H38").Value = 2367.6
$ws.Range("I38").Value = 1779.3334
$ws.Range("K38").Value = 5338.0002
$ws.Range("M38").Value = -4991.0002
$ws.Range("H61").Value = 131
$ws.Range("I61").Value = 67
$ws.Range("K61").Value = 201
$ws.Range("M61").Value = 14

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 52.454544
$ws.Range("I2").Value = 52.454544
$ws.Range("K2").Value = 52.454544
$ws.Range("M2").Value = 60.545456
$ws.Range("H113").Value = 3079
$ws.Range("I113").Value = 2743.5
$ws.Range("K113").Value = 2743.5
$ws.Range("M113").Value = -573.5
$ws.Range("H126").Value = 5416.7144
$ws.Range("I126").Value = 3148.4443
$ws.Range("K126").Value = 9445.332900000001
$ws.Range("M126").Value = -6975.332900000001

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H3").Value = 7232.5
$ws.Range("J3").Value = 10000
$ws.Range("L3").Value = 10000
$ws.Range("N3").Value = -10224
$ws.Range("H7").Value = 3190.4348
$ws.Range("I7").Value = 3208.6191
$ws.Range("K7").Value = 3208.6191
$ws.Range("M7").Value = -3096.6191
$ws.Range("H15").Value = 7232.5
$ws.Range("J15").Value = 10000
$ws.Range("L15").Value = 10000
$ws.Range("N15").Value = -10340
$ws.Range("H25").Value = 0
$ws.Range("J25").Value = 0
$ws.Range("L25").Value = 0
$ws.Range("H126").Value = 3190.4348
$ws.Range("I126").Value = 3208.6191
$ws.Range("K126").Value = 9625.8573
$ws.Range("M126").Value = -7155.8573
$ws.Range("H132").Value = 13184.737
$ws.Range("I132").Value = 13184.737
$ws.Range("K132").Value = 39554.211
$ws.Range("M132").Value = -37024.211
$ws.Range("N25").ClearContents()

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 2393
$ws.Range("I113").Value = 1549.5
$ws.Range("J113").Value = 3517.6667
$ws.Range("K113").Value = 4648.5
$ws.Range("L113").Value = 10553.0001
$ws.Range("M113").Value = -2478.5
$ws.Range("N113").Value = -14893.0001
$ws.Range("H135").Value = 96966.336
$ws.Range("J135").Value = 96966.336
$ws.Range("L135").Value = 96966.336
$ws.Range("N135").Value = -107106.336
